$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (F4 -> F2 suffix)
$ws.Range("B1").Value = "SPF_UNEMPF2"
$ws.Range("C1").Value = "SPF_INFF2"
$ws.Range("D1").Value = "SPF_COREINFF2"
$ws.Range("E1").Value = "SPF_gRPCEF2"

# Remove rows 46-68 (23 rows) so data ends at row 45
$ws.Range("A46:E68").EntireRow.Delete()
